$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 4414.517  # H62: 4483.357 -> 4414.517
$ws.Cells.Item(62, 9).Value = 3942.4211  # I62: 4023.2778 -> 3942.4211
$ws.Cells.Item(62, 11).Value = 3942.4211  # K62: 4023.2778 -> 3942.4211
$ws.Cells.Item(62, 13).Value = -3318.4211  # M62: -3399.2778 -> -3318.4211
$ws.Cells.Item(65, 8).Value = 4414.517  # H65: 4483.357 -> 4414.517
$ws.Cells.Item(65, 9).Value = 3942.4211  # I65: 4023.2778 -> 3942.4211
$ws.Cells.Item(65, 11).Value = 19712.1055  # K65: 20116.389 -> 19712.1055
$ws.Cells.Item(65, 13).Value = -16592.1055  # M65: -16996.389 -> -16592.1055
$ws.Cells.Item(106, 8).Value = 4754179  # H106: 4754188.5 -> 4754179
$ws.Cells.Item(106, 9).Value = 5373093.5  # I106: 5373104.5 -> 5373093.5
$ws.Cells.Item(106, 11).Value = 5373093.5  # K106: 5373104.5 -> 5373093.5
$ws.Cells.Item(106, 13).Value = -5372462.5  # M106: -5372473.5 -> -5372462.5
$ws.Cells.Item(132, 8).Value = 1541483.1  # H132: 1616091.1 -> 1541483.1
$ws.Cells.Item(132, 9).Value = 3121.14  # I132: 3313.9575 -> 3121.14
$ws.Cells.Item(132, 10).Value = 6669356.5  # J132: 6669460 -> 6669356.5
$ws.Cells.Item(132, 11).Value = 9363.42  # K132: 9941.872499999999 -> 9363.42
$ws.Cells.Item(132, 12).Value = 20008069.5  # L132: 20008380 -> 20008069.5
$ws.Cells.Item(132, 13).Value = -6833.42  # M132: -7411.872499999999 -> -6833.42
$ws.Cells.Item(132, 14).Value = -20013129.5  # N132: -20013440 -> -20013129.5
$ws.Cells.Item(137, 8).Value = 6290.217  # H137: 5977.6094 -> 6290.217
$ws.Cells.Item(137, 9).Value = 7510.4883  # I137: 7536.0933 -> 7510.4883
$ws.Cells.Item(137, 10).Value = 3203.647  # J137: 2786.4285 -> 3203.647
$ws.Cells.Item(137, 11).Value = 22531.4649  # K137: 22608.2799 -> 22531.4649
$ws.Cells.Item(137, 12).Value = 9610.940999999999  # L137: 8359.2855 -> 9610.940999999999
$ws.Cells.Item(137, 13).Value = -19981.4649  # M137: -20058.2799 -> -19981.4649
$ws.Cells.Item(137, 14).Value = -14710.941  # N137: -13459.2855 -> -14710.941
$ws.Cells.Item(138, 8).Value = 3827.473  # H138: 3902.081 -> 3827.473
$ws.Cells.Item(138, 9).Value = 1769.6072  # I138: 1856.5 -> 1769.6072
$ws.Cells.Item(138, 10).Value = 5080.087  # J138: 5010.104 -> 5080.087
$ws.Cells.Item(138, 11).Value = 5308.821599999999  # K138: 5569.5 -> 5308.821599999999
$ws.Cells.Item(138, 12).Value = 15240.261  # L138: 15030.312 -> 15240.261
$ws.Cells.Item(138, 13).Value = -168.8215999999993  # M138: -429.5 -> -168.8215999999993
$ws.Cells.Item(138, 14).Value = -25520.261  # N138: -25310.312 -> -25520.261
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 847.5  # H5: 850 -> 847.5
$ws.Cells.Item(5, 9).Value = 845  # I5: 0 -> 845
$ws.Cells.Item(5, 11).Value = 845  # K5: 0 -> 845
$ws.Cells.Item(5, 13).Value = -733  # M5: None -> -733
$ws.Cells.Item(61, 8).Value = 5365.8335  # H61: 5153.7354 -> 5365.8335
$ws.Cells.Item(61, 9).Value = 5075.4365  # I61: 4832.5967 -> 5075.4365
$ws.Cells.Item(61, 11).Value = 5075.4365  # K61: 4832.5967 -> 5075.4365
$ws.Cells.Item(61, 13).Value = -4863.4365  # M61: -4620.5967 -> -4863.4365
$ws.Cells.Item(80, 8).Value = 84974.875  # H80: 84898.5 -> 84974.875
$ws.Cells.Item(80, 10).Value = 84974.875  # J80: 84898.5 -> 84974.875
$ws.Cells.Item(80, 12).Value = 84974.875  # L80: 84898.5 -> 84974.875
$ws.Cells.Item(80, 14).Value = -86970.875  # N80: -86894.5 -> -86970.875
$ws.Cells.Item(83, 8).Value = 84974.875  # H83: 84898.5 -> 84974.875
$ws.Cells.Item(83, 10).Value = 84974.875  # J83: 84898.5 -> 84974.875
$ws.Cells.Item(83, 12).Value = 254924.625  # L83: 254695.5 -> 254924.625
$ws.Cells.Item(83, 14).Value = -264908.625  # N83: -264679.5 -> -264908.625
$ws.Cells.Item(132, 8).Value = 2020.303  # H132: 2065.3125 -> 2020.303
$ws.Cells.Item(132, 9).Value = 1210.0741  # I132: 1234.3077 -> 1210.0741
$ws.Cells.Item(132, 11).Value = 3630.2223  # K132: 3702.9231 -> 3630.2223
$ws.Cells.Item(132, 13).Value = -1100.2223  # M132: -1172.9231 -> -1100.2223
$ws.Cells.Item(136, 8).Value = 5365.8335  # H136: 5153.7354 -> 5365.8335
$ws.Cells.Item(136, 9).Value = 5075.4365  # I136: 4832.5967 -> 5075.4365
$ws.Cells.Item(136, 11).Value = 15226.3095  # K136: 14497.7901 -> 15226.3095
$ws.Cells.Item(136, 13).Value = -12676.3095  # M136: -11947.7901 -> -12676.3095
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 847.5  # H4: 850 -> 847.5
$ws.Cells.Item(4, 9).Value = 845  # I4: 0 -> 845
$ws.Cells.Item(4, 11).Value = 845  # K4: 0 -> 845
$ws.Cells.Item(4, 13).Value = -730  # M4: None -> -730
$ws.Cells.Item(38, 8).Value = 35000  # H38: 0 -> 35000
$ws.Cells.Item(38, 10).Value = 35000  # J38: 0 -> 35000
$ws.Cells.Item(38, 12).Value = 35000  # L38: 0 -> 35000
$ws.Cells.Item(38, 14).Value = -35832  # N38: None -> -35832
$ws.Cells.Item(134, 8).Value = 3250.4  # H134: 3363.6316 -> 3250.4
$ws.Cells.Item(134, 9).Value = 2143.7144  # I134: 2224.077 -> 2143.7144
$ws.Cells.Item(134, 11).Value = 6431.1432  # K134: 6672.231000000001 -> 6431.1432
$ws.Cells.Item(134, 13).Value = -3896.1432  # M134: -4137.231000000001 -> -3896.1432
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 6505.05  # H86: 6512.2 -> 6505.05
$ws.Cells.Item(86, 9).Value = 5282.8667  # I86: 5388.2856 -> 5282.8667
$ws.Cells.Item(86, 10).Value = 10171.6  # J86: 9134.666999999999 -> 10171.6
$ws.Cells.Item(86, 11).Value = 5282.8667  # K86: 5388.2856 -> 5282.8667
$ws.Cells.Item(86, 12).Value = 10171.6  # L86: 9134.666999999999 -> 10171.6
$ws.Cells.Item(86, 13).Value = -4159.8667  # M86: -4265.2856 -> -4159.8667
$ws.Cells.Item(86, 14).Value = -12417.6  # N86: -11380.667 -> -12417.6
$ws.Cells.Item(89, 8).Value = 6505.05  # H89: 6512.2 -> 6505.05
$ws.Cells.Item(89, 9).Value = 5282.8667  # I89: 5388.2856 -> 5282.8667
$ws.Cells.Item(89, 10).Value = 10171.6  # J89: 9134.666999999999 -> 10171.6
$ws.Cells.Item(89, 11).Value = 26414.3335  # K89: 26941.428 -> 26414.3335
$ws.Cells.Item(89, 12).Value = 50858  # L89: 45673.335 -> 50858
$ws.Cells.Item(89, 13).Value = -20798.3335  # M89: -21325.428 -> -20798.3335
$ws.Cells.Item(89, 14).Value = -62090  # N89: -56905.335 -> -62090
$ws.Cells.Item(134, 8).Value = 2595.8333  # H134: 2713.2942 -> 2595.8333
$ws.Cells.Item(134, 9).Value = 1294.4375  # I134: 1340.8 -> 1294.4375
$ws.Cells.Item(134, 11).Value = 3883.3125  # K134: 4022.4 -> 3883.3125
$ws.Cells.Item(134, 13).Value = -1348.3125  # M134: -1487.4 -> -1348.3125
$ws.Cells.Item(141, 8).Value = 107961.48  # H141: 110250.57 -> 107961.48
$ws.Cells.Item(141, 10).Value = 112277.29  # J141: 115261.63 -> 112277.29
$ws.Cells.Item(141, 12).Value = 112277.29  # L141: 115261.63 -> 112277.29
$ws.Cells.Item(141, 14).Value = -122637.29  # N141: -125621.63 -> -122637.29
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 106.181816  # H12: 111.818184 -> 106.181816
$ws.Cells.Item(12, 9).Value = 100.8  # I12: 116.5 -> 100.8
$ws.Cells.Item(12, 10).Value = 110.666664  # J12: 109.14286 -> 110.666664
$ws.Cells.Item(12, 11).Value = 302.4  # K12: 349.5 -> 302.4
$ws.Cells.Item(12, 12).Value = 331.999992  # L12: 327.42858 -> 331.999992
$ws.Cells.Item(12, 13).Value = -129.4  # M12: -176.5 -> -129.4
$ws.Cells.Item(12, 14).Value = -677.999992  # N12: -673.42858 -> -677.999992
$ws.Cells.Item(34, 8).Value = 2042.5217  # H34: 2033.0416 -> 2042.5217
$ws.Cells.Item(34, 9).Value = 1090.4  # I34: 999.4545000000001 -> 1090.4
$ws.Cells.Item(34, 10).Value = 2774.923  # J34: 2907.6155 -> 2774.923
$ws.Cells.Item(34, 11).Value = 3271.2  # K34: 2998.3635 -> 3271.2
$ws.Cells.Item(34, 12).Value = 8324.769  # L34: 8722.8465 -> 8324.769
$ws.Cells.Item(34, 13).Value = -3187.2  # M34: -2914.3635 -> -3187.2
$ws.Cells.Item(34, 14).Value = -8492.769  # N34: -8890.8465 -> -8492.769
$ws.Cells.Item(114, 8).Value = 2274.6667  # H114: 1997.4286 -> 2274.6667
$ws.Cells.Item(114, 10).Value = 3872  # J114: 4499 -> 3872
$ws.Cells.Item(114, 12).Value = 11616  # L114: 13497 -> 11616
$ws.Cells.Item(114, 14).Value = -18124  # N114: -20005 -> -18124
$ws.Cells.Item(131, 8).Value = 27030618  # H131: 20836582 -> 27030618
$ws.Cells.Item(131, 9).Value = 100008120  # I131: 200015730 -> 100008120
$ws.Cells.Item(131, 10).Value = 1914.8889  # J131: 1797.721 -> 1914.8889
$ws.Cells.Item(131, 11).Value = 300024360  # K131: 600047190 -> 300024360
$ws.Cells.Item(131, 12).Value = 5744.6667  # L131: 5393.163 -> 5744.6667
$ws.Cells.Item(131, 13).Value = -300019320  # M131: -600042150 -> -300019320
$ws.Cells.Item(131, 14).Value = -15824.6667  # N131: -15473.163 -> -15824.6667
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3599.5557  # H132: 3800 -> 3599.5557
$ws.Cells.Item(132, 9).Value = 3517.1765  # I132: 3720 -> 3517.1765
$ws.Cells.Item(132, 11).Value = 10551.5295  # K132: 11160 -> 10551.5295
$ws.Cells.Item(132, 13).Value = -8021.529500000001  # M132: -8630 -> -8021.529500000001
$ws.Cells.Item(139, 8).Value = 47115.1  # H139: 57200 -> 47115.1
$ws.Cells.Item(139, 10).Value = 47115.1  # J139: 57200 -> 47115.1
$ws.Cells.Item(139, 12).Value = 47115.1  # L139: 57200 -> 47115.1
$ws.Cells.Item(139, 14).Value = -57395.1  # N139: -67480 -> -57395.1
$ws.Cells.Item(141, 8).Value = 97320.164  # H141: 98464.39999999999 -> 97320.164
$ws.Cells.Item(141, 10).Value = 97320.164  # J141: 98464.39999999999 -> 97320.164
$ws.Cells.Item(141, 12).Value = 97320.164  # L141: 98464.39999999999 -> 97320.164
$ws.Cells.Item(141, 14).Value = -107680.164  # N141: -108824.4 -> -107680.164
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 874.1429000000001  # H22: 938.3333 -> 874.1429000000001
$ws.Cells.Item(22, 9).Value = 623.8  # I22: 626 -> 623.8
$ws.Cells.Item(22, 10).Value = 1500  # J22: 2500 -> 1500
$ws.Cells.Item(22, 11).Value = 623.8  # K22: 626 -> 623.8
$ws.Cells.Item(22, 12).Value = 1500  # L22: 2500 -> 1500
$ws.Cells.Item(22, 13).Value = -328.8  # M22: -331 -> -328.8
$ws.Cells.Item(22, 14).Value = -2090  # N22: -3090 -> -2090
$ws.Cells.Item(27, 8).Value = 874.1429000000001  # H27: 938.3333 -> 874.1429000000001
$ws.Cells.Item(27, 9).Value = 623.8  # I27: 626 -> 623.8
$ws.Cells.Item(27, 10).Value = 1500  # J27: 2500 -> 1500
$ws.Cells.Item(27, 11).Value = 623.8  # K27: 626 -> 623.8
$ws.Cells.Item(27, 12).Value = 1500  # L27: 2500 -> 1500
$ws.Cells.Item(27, 13).Value = -516.8  # M27: -519 -> -516.8
$ws.Cells.Item(27, 14).Value = -1714  # N27: -2714 -> -1714
$ws.Cells.Item(46, 8).Value = 2409.92  # H46: 2489.5 -> 2409.92
$ws.Cells.Item(46, 9).Value = 755.55554  # I46: 787.5 -> 755.55554
$ws.Cells.Item(46, 11).Value = 755.55554  # K46: 787.5 -> 755.55554
$ws.Cells.Item(46, 13).Value = -567.55554  # M46: -599.5 -> -567.55554
$ws.Cells.Item(93, 8).Value = 3258.7585  # H93: 3373.1482 -> 3258.7585
$ws.Cells.Item(93, 10).Value = 2531.4443  # J93: 2764.8572 -> 2531.4443
$ws.Cells.Item(93, 12).Value = 2531.4443  # L93: 2764.8572 -> 2531.4443
$ws.Cells.Item(93, 14).Value = -5027.4443  # N93: -5260.8572 -> -5027.4443
$ws.Cells.Item(97, 8).Value = 0  # H97: 30000 -> 0
$ws.Cells.Item(97, 10).Value = 0  # J97: 30000 -> 0
$ws.Cells.Item(97, 12).ClearContents()  # L97 was 30000
$ws.Cells.Item(97, 14).Value = 0  # N97: -31982 -> 0
$ws.Cells.Item(103, 8).Value = 0  # H103: 27500 -> 0
$ws.Cells.Item(103, 10).Value = 0  # J103: 27500 -> 0
$ws.Cells.Item(103, 12).ClearContents()  # L103 was 27500
$ws.Cells.Item(103, 14).Value = 0  # N103: -29844 -> 0
$ws.Cells.Item(106, 8).Value = 20000  # H106: 9272 -> 20000
$ws.Cells.Item(106, 10).Value = 20000  # J106: 9272 -> 20000
$ws.Cells.Item(106, 12).Value = 20000  # L106: 9272 -> 20000
$ws.Cells.Item(106, 14).Value = -22524  # N106: -11796 -> -22524
$ws.Cells.Item(118, 8).Value = 0  # H118: 25000 -> 0
$ws.Cells.Item(118, 9).Value = 0  # I118: 25000 -> 0
$ws.Cells.Item(118, 11).Value = 0  # K118: 25000 -> 0
$ws.Cells.Item(118, 13).ClearContents()  # M118 was -23343
$ws.Cells.Item(122, 8).Value = 5575.6924  # H122: 5802.4165 -> 5575.6924
$ws.Cells.Item(122, 9).Value = 5114.1055  # I122: 5379.8823 -> 5114.1055
$ws.Cells.Item(122, 11).Value = 15342.3165  # K122: 16139.6469 -> 15342.3165
$ws.Cells.Item(122, 13).Value = -12892.3165  # M122: -13689.6469 -> -12892.3165
$ws.Cells.Item(132, 8).Value = 308827.4  # H132: 321872.9 -> 308827.4
$ws.Cells.Item(132, 9).Value = 430101.7  # I132: 456031.6 -> 430101.7
$ws.Cells.Item(132, 11).Value = 1290305.1  # K132: 1368094.8 -> 1290305.1
$ws.Cells.Item(132, 13).Value = -1287775.1  # M132: -1365564.8 -> -1287775.1
$ws.Cells.Item(136, 8).Value = 4343  # H136: 4363.625 -> 4343
$ws.Cells.Item(136, 9).Value = 3341.4814  # I136: 3434.08 -> 3341.4814
$ws.Cells.Item(136, 10).Value = 5572.136  # J136: 5374 -> 5572.136
$ws.Cells.Item(136, 11).Value = 10024.4442  # K136: 10302.24 -> 10024.4442
$ws.Cells.Item(136, 12).Value = 16716.408  # L136: 16122 -> 16716.408
$ws.Cells.Item(136, 13).Value = -7474.4442  # M136: -7752.24 -> -7474.4442
$ws.Cells.Item(136, 14).Value = -21816.408  # N136: -21222 -> -21816.408
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 40000  # H54: 0 -> 40000
$ws.Cells.Item(54, 9).Value = 40000  # I54: 0 -> 40000
$ws.Cells.Item(54, 11).Value = 40000  # K54: 0 -> 40000
$ws.Cells.Item(54, 13).Value = -39480  # M54: None -> -39480
$ws.Cells.Item(132, 8).Value = 7967.0195  # H132: 8409.583000000001 -> 7967.0195
$ws.Cells.Item(132, 9).Value = 8859  # I132: 9493.736999999999 -> 8859
$ws.Cells.Item(132, 10).Value = 4309.9  # J132: 4289.8 -> 4309.9
$ws.Cells.Item(132, 11).Value = 26577  # K132: 28481.211 -> 26577
$ws.Cells.Item(132, 12).Value = 12929.7  # L132: 12869.4 -> 12929.7
$ws.Cells.Item(132, 13).Value = -24047  # M132: -25951.211 -> -24047
$ws.Cells.Item(132, 14).Value = -17989.7  # N132: -17929.4 -> -17989.7
$ws.Cells.Item(136, 8).Value = 316113.9  # H136: 309821.6 -> 316113.9
$ws.Cells.Item(136, 9).Value = 359742.6  # I136: 359756.53 -> 359742.6
$ws.Cells.Item(136, 10).Value = 3441.8333  # J136: 3078.4285 -> 3441.8333
$ws.Cells.Item(136, 11).Value = 1079227.8  # K136: 1079269.59 -> 1079227.8
$ws.Cells.Item(136, 12).Value = 10325.4999  # L136: 9235.2855 -> 10325.4999
$ws.Cells.Item(136, 13).Value = -1076677.8  # M136: -1076719.59 -> -1076677.8
$ws.Cells.Item(136, 14).Value = -15425.4999  # N136: -14335.2855 -> -15425.4999
